$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RUNMANAGER
$ws2 = $wb.Worksheets.Item(2)   # DATA

# ---------------------------------------------------------------------------
# RUNMANAGER sheet
# ---------------------------------------------------------------------------

# multiplyLoginLogoutTest row is now executed
$ws1.Range("C4").Value = "yes"

# The old "multiplySubMenuOptionsTest" row is removed; remaining rows shift up
$ws1.Rows.Item(5).Delete() | Out-Null

# New test row describing the multiply rewards test
$ws1.Range("A6").Value = "multiplyRewardsTest"
$ws1.Range("B6").Value = "To check multiply rewards test"
$ws1.Range("C6").Value = "no"
$ws1.Range("D6").Value = "'1"
$ws1.Range("E6").Value = "'1"

# ---------------------------------------------------------------------------
# DATA sheet
# ---------------------------------------------------------------------------

# multiplyLoginLogoutTest data row now executes
$ws2.Range("B7").Value = "yes"
# multiplySubMenuOptionsTest data row no longer executes
$ws2.Range("B8").Value = "no"

# New data row for multiplyRewardsTest
$ws2.Range("A10").Value = "multiplyRewardsTest"
$ws2.Range("B10").Value = "yes"
$ws2.Range("C10").Value = "chrome"
$ws2.Range("D10").Value = "'"
$ws2.Range("E10").Value = "'"

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------

# Set DATA's lingering selection first (it is no longer the active tab)
$ws2.Range("E8").Select() | Out-Null

# RUNMANAGER becomes the active tab with its own new selection
$ws1.Activate() | Out-Null
$ws1.Range("C6").Select() | Out-Null
